$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.254.86"
$ws.Range("E2").Value = "  +4.02%  "
$ws.Range("D3").Value = "2.341.64"
$ws.Range("E3").Value = "  +2.54%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "545.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.94%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("D9").Value = "2.337.91"
$ws.Range("E9").Value = "  +2.43%  "
$ws.Range("E10").Value = "  +1.91%  "
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.03%  "
$ws.Range("D15").Value = "2.757.48"
$ws.Range("E15").Value = "  +2.46%  "
$ws.Range("D16").Value = "60.184.92"
$ws.Range("E16").Value = "  +3.99%  "
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").Value = "2.335.15"
$ws.Range("E18").Value = "  +2.48%  "
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "313.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.04%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("E25").Value = "  +3.20%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("E28").Value = "  +8.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("E31").Value = "  +13.49%  "
$ws.Range("D32").Value = "0.0₃0730"
$ws.Range("E32").Value = "  +2.41%  "
$ws.Range("E33").Value = "  +4.24%  "
$ws.Range("E34").Value = "  +13.81%  "
$ws.Range("E35").Value = "  +1.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  +7.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "322.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "141.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("E44").Value = "  +1.95%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.38%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0945"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.560"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0213"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("D51").Value = "0.0₆0210"
$ws.Range("E51").Value = "  +18.00%  "
